$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update period label for row 2: 2014/12 -> 2015/03
$ws.Range("C2").Value = "2015/03  (IFRS연결)"

# Row 2
$ws.Range("D2").Value = 2118
$ws.Range("E2").Value = 62
$ws.Range("F2").Value = 62
$ws.Range("G2").Value = -27
$ws.Range("H2").Value = -57
$ws.Range("I2").Value = -110
$ws.Range("J2").Value = 53
$ws.Range("K2").Value = 4308
$ws.Range("L2").Value = 2474
$ws.Range("M2").Value = 1835
$ws.Range("N2").Value = 1734
$ws.Range("O2").Value = 100
$ws.Range("P2").Value = 459
$ws.Range("Q2").Value = 185
$ws.Range("R2").Value = -204
$ws.Range("S2").Value = 112
$ws.Range("T2").Value = 139
$ws.Range("U2").Value = 46
$ws.Range("V2").Value = 1910
$ws.Range("W2").Value = 2.95
$ws.Range("X2").Value = -2.67
$ws.Range("AA2").Value = 134.83
$ws.Range("AB2").Value = 162.82
$ws.Range("AC2").Value = -598
$ws.Range("AD2").Value = -56.8
$ws.Range("AE2").Value = 9934
$ws.Range("AF2").Value = 3.42
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 17895256

# Row 3
$ws.Range("D3").Value = 1863
$ws.Range("E3").Value = 155
$ws.Range("F3").Value = 155
$ws.Range("G3").Value = 115
$ws.Range("H3").Value = 83
$ws.Range("I3").Value = 36
$ws.Range("J3").Value = 47
$ws.Range("K3").Value = 4316
$ws.Range("L3").Value = 2206
$ws.Range("M3").Value = 2110
$ws.Range("N3").Value = 1949
$ws.Range("O3").Value = 161
$ws.Range("P3").Value = 487
$ws.Range("Q3").Value = 218
$ws.Range("R3").Value = -219
$ws.Range("S3").Value = -51
$ws.Range("T3").Value = 77
$ws.Range("U3").Value = 141
$ws.Range("V3").Value = 1641
$ws.Range("W3").Value = 8.33
$ws.Range("X3").Value = 4.46
$ws.Range("Y3").Value = 1.96
$ws.Range("Z3").Value = 1.93
$ws.Range("AA3").Value = 104.56
$ws.Range("AB3").Value = 192.89
$ws.Range("AC3").Value = 189
$ws.Range("AD3").Value = 259.28
$ws.Range("AE3").Value = 10474
$ws.Range("AF3").Value = 4.68
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 19045984

# Row 4
$ws.Range("D4").Value = 2616
$ws.Range("E4").Value = 225
$ws.Range("F4").Value = 225
$ws.Range("G4").Value = 162
$ws.Range("H4").Value = 119
$ws.Range("I4").Value = 58
$ws.Range("J4").Value = 61
$ws.Range("K4").Value = 4212
$ws.Range("L4").Value = 2139
$ws.Range("M4").Value = 2073
$ws.Range("N4").Value = 1895
$ws.Range("O4").Value = 177
$ws.Range("P4").Value = 487
$ws.Range("Q4").Value = 270
$ws.Range("R4").Value = -106
$ws.Range("S4").Value = -123
$ws.Range("T4").Value = 27
$ws.Range("U4").Value = 242
$ws.Range("V4").Value = 1550
$ws.Range("W4").Value = 8.6
$ws.Range("X4").Value = 4.55
$ws.Range("Y4").Value = 3.04
$ws.Range("Z4").Value = 2.79
$ws.Range("AA4").Value = 103.2
$ws.Range("AB4").Value = 183.42
$ws.Range("AC4").Value = 300
$ws.Range("AD4").Value = 138.83
$ws.Range("AE4").Value = 10185
$ws.Range("AF4").Value = 4.09
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 0.24
$ws.Range("AI4").Value = 32.01
$ws.Range("AJ4").Value = 19045984

# Row 5
$ws.Range("D5").Value = 2698
$ws.Range("E5").Value = 244
$ws.Range("F5").Value = 244
$ws.Range("G5").Value = 182
$ws.Range("H5").Value = 129
$ws.Range("I5").Value = 63
$ws.Range("J5").Value = 67
$ws.Range("K5").Value = 4364
$ws.Range("L5").Value = 2254
$ws.Range("M5").Value = 2110
$ws.Range("N5").Value = 1919
$ws.Range("O5").Value = 190
$ws.Range("P5").Value = 487
$ws.Range("Q5").Value = 300
$ws.Range("R5").Value = -130
$ws.Range("S5").Value = -58
$ws.Range("T5").Value = 38
$ws.Range("U5").Value = 262
$ws.Range("V5").Value = 1544
$ws.Range("W5").Value = 9.06
$ws.Range("X5").Value = 4.8
$ws.Range("Y5").Value = 3.28
$ws.Range("Z5").Value = 3.02
$ws.Range("AA5").Value = 106.83
$ws.Range("AB5").Value = 191.71
$ws.Range("AC5").Value = 321
$ws.Range("AD5").Value = 119.94
$ws.Range("AE5").Value = 10314
$ws.Range("AF5").Value = 3.74
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 19045984

# Row 6
$ws.Range("D6").Value = 3000
$ws.Range("E6").Value = 167
$ws.Range("F6").Value = 167
$ws.Range("G6").Value = 81
$ws.Range("H6").Value = 33
$ws.Range("I6").Value = -49
$ws.Range("K6").Value = 4323
$ws.Range("L6").Value = 2224
$ws.Range("M6").Value = 2099
$ws.Range("N6").Value = 1875
$ws.Range("P6").Value = 488
$ws.Range("Q6").Value = 221
$ws.Range("R6").Value = -114
$ws.Range("S6").Value = -115
$ws.Range("T6").Value = 15
$ws.Range("U6").Value = 206
$ws.Range("V6").Value = 1440
$ws.Range("W6").Value = 5.57
$ws.Range("X6").Value = 1.08
$ws.Range("Y6").Value = -2.58
$ws.Range("Z6").Value = 0.75
$ws.Range("AA6").Value = 105.95
$ws.Range("AB6").Value = 184.36
$ws.Range("AC6").Value = -251
$ws.Range("AD6").Value = -109.55
$ws.Range("AE6").Value = 10053
$ws.Range("AF6").Value = 2.74
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 19085664

# Row 7
$ws.Range("D7").Value = 3280
$ws.Range("E7").Value = 315
$ws.Range("G7").Value = 264
$ws.Range("H7").Value = 197
$ws.Range("I7").Value = 117
$ws.Range("K7").Value = 4485
$ws.Range("L7").Value = 2193
$ws.Range("M7").Value = 2292
$ws.Range("N7").Value = 2003
$ws.Range("P7").Value = 488
$ws.Range("Q7").Value = 168
$ws.Range("R7").Value = -50
$ws.Range("S7").Value = -73
$ws.Range("T7").Value = 9
$ws.Range("W7").Value = 9.6
$ws.Range("X7").Value = 6.01
$ws.Range("Y7").Value = 6.03
$ws.Range("Z7").Value = 4.47
$ws.Range("AA7").Value = 95.68000000000001
$ws.Range("AC7").Value = 599
$ws.Range("AD7").Value = 40.15
$ws.Range("AE7").Value = 10740
$ws.Range("AF7").Value = 2.24
$ws.Range("AG7").Value = 100
$ws.Range("AH7").Value = 0.42
$ws.Range("AI7").Value = 16.31

# Row 8
$ws.Range("D8").Value = 3572
$ws.Range("E8").Value = 348
$ws.Range("G8").Value = 304
$ws.Range("H8").Value = 225
$ws.Range("I8").Value = 139
$ws.Range("K8").Value = 4893
$ws.Range("L8").Value = 2396
$ws.Range("M8").Value = 2497
$ws.Range("N8").Value = 2208
$ws.Range("P8").Value = 488
$ws.Range("Q8").Value = 388
$ws.Range("R8").Value = -61
$ws.Range("S8").Value = 130
$ws.Range("T8").Value = 18
$ws.Range("W8").Value = 9.74
$ws.Range("X8").Value = 6.3
$ws.Range("Y8").Value = 6.6
$ws.Range("Z8").Value = 4.8
$ws.Range("AA8").Value = 95.95999999999999
$ws.Range("AC8").Value = 712
$ws.Range("AD8").Value = 33.79
$ws.Range("AE8").Value = 11839
$ws.Range("AF8").Value = 2.03
$ws.Range("AG8").Value = 100
$ws.Range("AH8").Value = 0.42
$ws.Range("AI8").Value = 13.73

# Row 9
$ws.Range("D9").Value = 3985
$ws.Range("E9").Value = 409
$ws.Range("G9").Value = 362
$ws.Range("H9").Value = 268
$ws.Range("I9").Value = 166
$ws.Range("K9").Value = 5211
$ws.Range("L9").Value = 2465
$ws.Range("M9").Value = 2746
$ws.Range("N9").Value = 2457
$ws.Range("P9").Value = 488
$ws.Range("Q9").Value = 386
$ws.Range("R9").Value = -64
$ws.Range("S9").Value = -20
$ws.Range("T9").Value = 21
$ws.Range("W9").Value = 10.26
$ws.Range("X9").Value = 6.73
$ws.Range("Y9").Value = 7.12
$ws.Range("Z9").Value = 5.3
$ws.Range("AA9").Value = 89.77
$ws.Range("AC9").Value = 850
$ws.Range("AD9").Value = 28.3
$ws.Range("AE9").Value = 13174
$ws.Range("AF9").Value = 1.83
$ws.Range("AG9").Value = 100
$ws.Range("AH9").Value = 0.42
$ws.Range("AI9").Value = 11.5

# Cells removed entirely (shift-left style clears) per diff
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("U9").ClearContents()